$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 20 / Row 21: coin swap (Chainlink <-> ShibaInu) ---
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "0.0₃0693"
$ws.Range("E20").Value = "  -2.68%  "

$ws.Range("B21").Value = "Chainlink"
$ws.Range("C21").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D21").Value = "'7.34"
$ws.Range("E21").Value = "  -2.95%  "

# --- Remaining rows: Price (D) and/or Volume(1h) (E) updates ---
$ws.Range("D2").Value = "28.473.36"
$ws.Range("E2").Value = "  -0.17%  "
$ws.Range("D3").Value = "1.564.85"
$ws.Range("E3").Value = "  -2.27%  "
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").Value = "'211.72"
$ws.Range("E5").Value = "  -1.59%  "
$ws.Range("E6").Value = "  -1.29%  "
$ws.Range("E7").Value = "  +0.29%  "
$ws.Range("D8").Value = "'46.15"
$ws.Range("E8").Value = "  +4.59%  "
$ws.Range("D9").Value = "'23.98"
$ws.Range("E9").Value = "  -0.69%  "
$ws.Range("E10").Value = "  -1.89%  "
$ws.Range("E11").Value = "  -1.67%  "
$ws.Range("D12").Value = "'0.0886"
$ws.Range("E12").Value = "  -0.43%  "
$ws.Range("D13").Value = "1.789.75"
$ws.Range("E13").Value = "  -2.08%  "
$ws.Range("D14").Value = "1.566.96"
$ws.Range("E14").Value = "  -2.13%  "
$ws.Range("E15").Value = "  -2.93%  "
$ws.Range("D16").Value = "28.489.15"
$ws.Range("E16").Value = "  -0.10%  "
$ws.Range("E17").Value = "  -3.62%  "
$ws.Range("D18").Value = "'62.21"
$ws.Range("E18").Value = "  -1.92%  "
$ws.Range("D19").Value = "'228.47"
$ws.Range("E19").Value = "  -2.08%  "
$ws.Range("E22").Value = "  +0.16%  "
$ws.Range("E23").Value = "  -6.06%  "
$ws.Range("D24").Value = "'9.12"
$ws.Range("E24").Value = "  -3.37%  "
$ws.Range("E25").Value = "  +6.33%  "
$ws.Range("D26").Value = "'150.86"
$ws.Range("E26").Value = "  -1.30%  "
$ws.Range("E27").Value = "  -2.29%  "
$ws.Range("E28").Value = "  -3.01%  "
$ws.Range("E29").Value = "  -4.26%  "
$ws.Range("E30").Value = "  +0.26%  "
$ws.Range("E31").Value = "  -2.15%  "
$ws.Range("D32").Value = "'1.11"
$ws.Range("E32").Value = "  -3.92%  "
$ws.Range("E33").Value = "  -1.53%  "
$ws.Range("E34").Value = "  -3.13%  "
$ws.Range("D35").Value = "1.390.05"
$ws.Range("E35").Value = "  -2.41%  "
$ws.Range("D36").Value = "'1.05"
$ws.Range("E36").Value = "  -1.39%  "
$ws.Range("E37").Value = "  -3.81%  "
$ws.Range("D38").Value = "'2.36"
$ws.Range("E38").Value = "  +1.13%  "
$ws.Range("E39").Value = "  +2.08%  "
$ws.Range("E40").Value = "  -1.36%  "
$ws.Range("D41").Value = "'0.535"
$ws.Range("E41").Value = "  -1.98%  "
$ws.Range("E42").Value = "  +0.22%  "
$ws.Range("D43").Value = "'1.88"
$ws.Range("E43").Value = "  +2.25%  "
$ws.Range("E44").Value = "  -4.41%  "
$ws.Range("E45").Value = "  -4.53%  "
$ws.Range("E46").Value = "  -0.32%  "
$ws.Range("D47").Value = "'62.69"
$ws.Range("E47").Value = "  -3.66%  "
$ws.Range("D48").Value = "1.701.57"
$ws.Range("E48").Value = "  -2.10%  "
$ws.Range("D49").Value = "'86.00"
$ws.Range("E49").Value = "  -1.96%  "
$ws.Range("E50").Value = "  -4.34%  "
$ws.Range("E51").Value = "  -0.44%  "
